# Apply "cetak susulan dp2nt16 dp3n31 1nt9" edits:
# Update mail-merge result fields (NO, NAMA, SEPATU, TOPI) for both
# label cells in the document.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Left label (C1 -> C73)
Replace-Text "C1" "C73"
Replace-Text "EFRAIN" "IDHO YUDHA F"
Replace-Text "40" "41"
Replace-Text "54" "58"

# Right label (C2 -> C74)
Replace-Text "C2" "C74"
Replace-Text "ANDY SETYO PRASONGKO" "M. ALSY SYARIFUDIN AF"
Replace-Text "42" "40"
Replace-Text "56" "54"

$d.Save()
